# Restructure the roster sheet:
#  - drop the old "Unnamed: 0" index column (B), using its values to
#    replace the row-index column (A)
#  - shift lastname/firstname/middlename/grade/email/password/sent left
#  - reset the "sent" flag to 0 for every row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the "Unnamed: 0" values (col B, rows 2-14) into col A, replacing
#    the old row index.
$ws.Range("A2:A14").Value2 = $ws.Range("B2:B14").Value2

# 2) Delete the now-redundant column B; everything to the right (lastname..sent)
#    shifts left by one column.
$ws.Range("B1").EntireColumn.Delete()

# 3) Reset the "sent" column (now H) to 0 for every data row.
$ws.Range("H2:H14").Value2 = 0
